$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New block header label (row 14)
$ws.Range("A14").Value = "Aussc hwing?"

# New 4x4 matrix block (rows 15-18)
$ws.Range("A15").Value = 0.01
$ws.Range("B15").Value = 0
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 0

$ws.Range("A16").Value = 0
$ws.Range("B16").Value = 10
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 0

$ws.Range("A17").Value = 0
$ws.Range("B17").Value = 0
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 0

$ws.Range("A18").Value = 0
$ws.Range("B18").Value = 0
$ws.Range("C18").Value = 0
$ws.Range("D18").Value = 0.0001

# Match the updated selection/active cell from the edited file
$ws.Range("A15:D18").Select()

# Page setup: paper size + orientation (matches added pageSetup element)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
